$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Replace the sentence "It doesn't matter what data type the CustomList
#    is." with the more specific wording, split across several runs (as the
#    author produced it through a handful of separate edits) while keeping
#    the original Georgia run formatting.
# ---------------------------------------------------------------------------

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*CustomList is.*") {
        $target = $cand
        break
    }
}

$pStart = $target.Range.Start
$pEnd = $target.Range.End - 1
$replaceRange = $d.Range($pStart, $pEnd)

$newRunsXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/></w:rPr><w:t>The</w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/></w:rPr><w:t xml:space="preserve"> data type </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/></w:rPr><w:t xml:space="preserve">for </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/></w:rPr><w:t xml:space="preserve">CustomList </w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/></w:rPr><w:t>can be any data type</w:t></w:r>' + `
'<w:r><w:rPr><w:rFonts w:ascii="Georgia" w:hAnsi="Georgia"/></w:rPr><w:t>.</w:t></w:r>' + `
'</w:p>' + `
'</w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

$replaceRange.InsertXML($newRunsXml)

# Re-resolve the paragraph (content shifted) for the bookmark move below.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*can be any data type.*") {
        $target = $cand
        break
    }
}

# ---------------------------------------------------------------------------
# 2. Move the hidden "_GoBack" bookmark from the end of the document (after
#    the final ")" run) to the end of the paragraph we just edited - this is
#    the usual Word behaviour of recording the most recent edit location.
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}

# Adding a bookmark with a collapsed range sitting exactly on the last
# character position of a paragraph (i.e. immediately before the paragraph
# mark) is mishandled by this host, so a one-character placeholder is typed,
# the bookmark is anchored just before it, and the placeholder is removed
# again - leaving the bookmark correctly collapsed right after the final
# "." run and before the paragraph mark.
$endPos = $target.Range.End - 1
$placeholderRange = $d.Range($endPos, $endPos)
$placeholderRange.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$delRange = $d.Range($endPos, $endPos + 1)
$delRange.Text = ""
